$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Relevant Courses" line: add "Computer Vision" and abbreviate
#    "Statistical Machine Learning" -> "Computer Vision, Statistical ML"
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Statistical Machine Learning, ", $false, $false, $false, $false, $false, $true, 1, $false, "Computer Vision, Statistical ML, ", 2)

# ------------------------------------------------------------------
# 2) Insert a new, collapsed "_GoBack" bookmark right after the
#    "2021" at the end of the M.S. education paragraph.
#    (Collapsed ranges landing exactly on a paragraph boundary are
#    mishandled by this runtime, so we briefly insert a one
#    character placeholder, anchor the bookmark to it, then delete
#    the placeholder text again -- the bookmark itself survives.)
# ------------------------------------------------------------------
$eduPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*2021*") {
        $eduPara = $p
    }
}
$endRng = $eduPara.Range
$endRng.MoveEnd(1, -1)
$endRng.Collapse(0)
$endRng.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $endRng)
$placeholder = $d.Range($endRng.Start, $endRng.End)
$placeholder.Text = ""

# ------------------------------------------------------------------
# 3) Merge "train" + (old "_GoBack" bookmark) + "ing progress..."
#    into a single run reading "training progress..." and drop the
#    stale bookmark that used to mark the editing cursor there.
# ------------------------------------------------------------------
$trainRng = $d.Content
$trainRng.Find.Execute("training progress and error rates of correlation for finished trainings")
$trainRng.Text = "ZZPLACEHOLDERZZtraining progress and error rates of correlation for finished trainings"

$trainRng2 = $d.Content
$trainRng2.Find.Execute("ZZPLACEHOLDERZZtraining progress and error rates of correlation for finished trainings", $false, $false, $false, $false, $false, $true, 1, $false, "training progress and error rates of correlation for finished trainings", 2)
